# CDS_SPARSE_2020_2021.xlsx merge-debug edit
# 1. Fix two swapped 0/1 cells on "High School Units" (sheet3: C7/D7).
# 2. Add a new "Admission_General" worksheet at the end of the workbook,
#    populate it with the admissions table, and leave selections/active
#    tab matching the target state.

$wb = $excel.ActiveWorkbook

# --- 1. High School Units: swap C7 (1 -> 0) and D7 (0 -> 1) ---------------
$hs = $wb.Worksheets.Item("High School Units")
$hs.Cells.Item(7, 3).Value = 0
$hs.Cells.Item(7, 4).Value = 1

# --- 2. Add Admission_General as the last sheet ----------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$admission = $wb.Worksheets.Add($null, $lastSheet)
$admission.Name = "Admission_General"

$data = @(
        @("Value","freshman","women","men","applied","admitted","enrolled","full-time","part-time","waiting-list-policy","offered-list","accepted-list","admitted-list","list-rank","completion-requirement","college-preparatory-program"),
        @(3338,1,0,1,1,0,0,0,0,0,0,0,0,0,0,0),
        @(1038,1,1,0,1,0,0,0,0,0,0,0,0,0,0,0),
        @(2502,1,0,1,0,1,0,0,0,0,0,0,0,0,0,0),
        @(851,1,1,0,0,1,0,0,0,0,0,0,0,0,0,0),
        @(401,1,0,1,0,0,1,1,0,0,0,0,0,0,0,0),
        @(1,1,0,1,0,0,1,0,1,0,0,0,0,0,0,0),
        @(145,1,1,0,0,0,1,1,0,0,0,0,0,0,0,0),
        @(0,1,1,0,0,0,1,0,1,0,0,0,0,0,0,0),
        @("Yes",0,0,0,0,0,0,0,0,1,0,0,0,0,0,0),
        @(119,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0),
        @(43,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0),
        @(14,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0),
        @("No",0,0,0,0,0,0,0,0,0,0,0,0,1,0,0),
        @("High school diploma is required and GED is not accepted",0,0,0,0,0,0,0,0,0,0,0,0,0,1,0),
        @("Require",0,0,0,0,0,0,0,0,0,0,0,0,0,0,1)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $admission.Cells.Item($r + 1, $c + 1).Value = $rowValues[$c]
    }
}

# --- 3. Selections: new sheet shows G18 selected, but "High School Units"
#        stays the active tab with E7 selected (matches target sheetViews).
$admission.Activate() | Out-Null
$admission.Range("G18").Select() | Out-Null

$hs.Activate() | Out-Null
$hs.Range("E7").Select() | Out-Null
